# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" colours (only used by the notes
#                            master relationship)
#   ppt/theme/theme2.xml -> "Integral" colours (the theme actually applied
#                            to the slide master / the whole deck)
#
# The target edit swaps the two themes' colour schemes (their font/format
# schemes are byte-identical, so swapping just the 12 theme colours
# reproduces the effective change). Re-colour the live/applied theme so it
# carries the "Office Theme" palette instead of "Integral".

$p = $ppt.ActivePresentation
$scheme = $p.SlideMaster.Theme.ThemeColorScheme

# MsoThemeColorSchemeIndex (1-12) -> target "Office Theme" colour, written
# through .RGB as VBA's 0x00BBGGRR long (R | G<<8 | B<<16).
$scheme.Colors(1).RGB  = 0           # dk1      #000000
$scheme.Colors(2).RGB  = 16777215    # lt1      #FFFFFF
$scheme.Colors(3).RGB  = 6968388     # dk2      #44546A
$scheme.Colors(4).RGB  = 15132391    # lt2      #E7E6E6
$scheme.Colors(5).RGB  = 13998939    # accent1  #5B9BD5
$scheme.Colors(6).RGB  = 3243501     # accent2  #ED7D31
$scheme.Colors(7).RGB  = 10855845    # accent3  #A5A5A5
$scheme.Colors(8).RGB  = 49407       # accent4  #FFC000
$scheme.Colors(9).RGB  = 12874308    # accent5  #4472C4
$scheme.Colors(10).RGB = 4697456     # accent6  #70AD47
$scheme.Colors(11).RGB = 12673797    # hlink    #0563C1
$scheme.Colors(12).RGB = 7491477     # folHlink #954F72
